# Applies the "added output check file" edit to the scenarios_as_columns sheet:
#  - row 209 J:AC divided by an extra factor of 5
#  - row 210 J:AC bumped from 0.25 to 0.3
#  - row 212 J:AC switched from a flat 0.5 literal to the "0.00002 / 14" formula
#  - row 214 J:AC bumped from 0 to 0.15
#  - row 220 J:AC switched from a flat literal to the "0.00000012*100" formula
#  - row 226 J:AC dropped from 0.2 to 0.1
#  - row 228 J:AC formula changed from "2*0.1" to "2*0.075"
#  - a brand-new parameter row (smax_AA) appended as row 230

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 209: J:AC now divide the existing rate by an additional factor of 5 ---
$ws.Range("J209:AC209").Formula = "=0.00002 / 14 / 5"

# --- Row 210: J:AC literal bump 0.25 -> 0.3 ---
$ws.Range("J210:AC210").Value = 0.3

# --- Row 212: J:AC go from a flat literal 0.5 to a real formula, matching row 209's format ---
$ws.Range("H209").Copy()
$ws.Range("J212:AC212").PasteSpecial(-4122)
$ws.Range("J212:AC212").Formula = "=0.00002 / 14"

# --- Row 214: J:AC literal bump 0 -> 0.15 ---
$ws.Range("J214:AC214").Value = 0.15

# --- Row 220: J:AC go from a flat literal to a real formula ---
$ws.Range("J220:AC220").Formula = "=0.00000012*100"

# --- Row 226: J:AC literal drop 0.2 -> 0.1 ---
$ws.Range("J226:AC226").Value = 0.1

# --- Row 228: J:AC formula changes from "2*0.1" to "2*0.075" (AE:AO stays "2*0.1") ---
$ws.Range("J228:AC228").Formula = "=2*0.075"

# --- New row 230: smax_AA parameter, cloned from row 229's layout/format then edited ---
$ws.Range("A229:AO229").Copy($ws.Range("A230:AO230"))
$ws.Range("G230:H230").NumberFormat = "0.00E+00"
$ws.Range("J230:AC230").NumberFormat = "0.00E+00"
$ws.Range("AE230:AO230").NumberFormat = "0.00E+00"
$ws.Range("I230").Clear()
$ws.Range("AD230").Clear()

$ws.Cells.Item(230, 1).Value = "smax_AA"
$ws.Cells.Item(230, 2).Value = "parameter"
$ws.Cells.Item(230, 3).Value = "root_cynaps"
$ws.Cells.Item(230, 4).Value = "roots"
$ws.Cells.Item(230, 5).Value = "Maximal rate of amino acid synthesis in the root segment"
$ws.Cells.Item(230, 6).Value = "m"

$ws.Range("G230:H230").Value = 0.00001
$ws.Range("J230:AC230").Formula = "=0.00001*10"
$ws.Range("AE230:AO230").Formula = "=0.00001"

# Leave the selection where the author last clicked, so the saved view matches.
$ws.Range("L228").Select()
